$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6495413184165955
$ws.Range("B1").Value = 2.159087896347046
$ws.Range("D1").Value = 1.054647922515869
$ws.Range("E1").Value = 1.162374138832092
